$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7287194209349384
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 0.7127328510149897
$ws.Range("E2").Value = 0.4998867070740569
$ws.Range("G2").Value = 3.594575437922795
